# Add two new columns, I ("I0") and J ("IF"), to the right of the existing
# H ("IP") column on the only worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# The existing headers (B1:H1) share one cell style (bold font, thin box
# border, centered/top-aligned) stored as style index 1. Clone that exact
# style onto the two new header cells by copying formats from H1, instead
# of rebuilding it property-by-property (which would mint a new style).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (rows 2-15) -------------------------------------------------
$data = @{
    2  = @(1, 4)
    3  = @(1, 6)
    4  = @(1, 5)
    5  = @(5, 7)
    6  = @(1, 5)
    7  = @(1, 4)
    8  = @(6, 9)
    9  = @(1, 5)
    10 = @(1, 4)
    11 = @(1, 5)
    12 = @(1, 4)
    13 = @(1, 5)
    14 = @(1, 6)
    15 = @(1, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
